$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J (|S*|/n) column, bold font
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# Row 14: Average of SW(S*)/SW(OPT)
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108

# Row 15: Average of SC(S*)/SC(OPT)
$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B15").Font.Bold = $true
$ws.Range("B15").Font.Size = 12
$ws.Range("B15").VerticalAlignment = -4108

# Row 16: Worst of SW(S*)/SW(OPT)
$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B16").Font.Bold = $true
$ws.Range("B16").Font.Size = 12
$ws.Range("B16").VerticalAlignment = -4108

# Row 17: Worst of SC(S*)/SC(OPT)
$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"
$ws.Range("B17").Font.Bold = $true
$ws.Range("B17").Font.Size = 12
$ws.Range("B17").VerticalAlignment = -4108

# Row heights for the new summary rows (15.6 points, as in the authored file)
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Select B17 (matches the saved cursor position in the authored file)
$ws.Range("B17").Select() | Out-Null

# Page setup matching the authored file
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
